$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1, "55555", "Deep Mondal", 34, "M", "Confirmed"),
    @(2, "55555", "Qwe uiio", 78, "F", "Confirmed"),
    @(3, "55555", "rttyy vvbbb", 98, "F", "Confirmed")
)

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]

    # Train Number looks like a plain number ("55555"); writing it straight
    # would let Excel infer a numeric cell (or, with an apostrophe, tag the
    # cell with a quote-prefix style). Route it through a text formula and
    # collapse that to a literal value so the stored cell is a clean shared
    # string with no extra formatting applied.
    $cellB = $ws.Cells.Item($rowIndex, 2)
    $cellB.Formula = '=TEXT(' + $row[1] + ',"0")'
    $cellB.Copy()
    $cellB.PasteSpecial(-4163)
    $excel.CutCopyMode = $false

    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    $ws.Cells.Item($rowIndex, 6).Value = $row[5]
    $rowIndex++
}
